$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 190 (shifts the existing rows 190-218 down to 191-219)
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A190").Value = 3
$ws.Range("B190").Value = "Femacal de La Calera"
$ws.Range("C190").Value = "Coquimbo"
$ws.Range("D190").Value = 44474
$ws.Range("E190").Value = 5
$ws.Range("F190").Value = 100112031
$ws.Range("G190").Value = "Poroto verde"
$ws.Range("H190").Value = "Magnum"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 38
$ws.Range("K190").Value = 35000
$ws.Range("L190").Value = 35000
$ws.Range("M190").Value = 35000
$ws.Range("N190").Value = "$/malla 25 kilos"
$ws.Range("O190").Value = "Región de Arica y Parinacota"
$ws.Range("P190").Value = 1400
$ws.Range("Q190").Value = 25
$ws.Range("R190").Value = "Hortaliza"
